$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3, 24.04999923706055, 28, 28.6299991607666, 22.65999984741211),
    @(4, 17.59000015258789, 22.1200008392334, 22.71999931335449, 16.71999931335449),
    @(5, 22.22999954223633, 20.42000007629395, 23.25, 18.95000076293945),
    @(6, 20.11000061035156, 24.71999931335449, 24.95000076293945, 19.70000076293945),
    @(7, 30.43000030517578, 35.59999847412109, 37.77000045776367, 30.30999946594238),
    @(8, 44.20000076293945, 40, 47, 38.90000152587891),
    @(9, 45.25, 52.54999923706055, 56.84999847412109, 44.45000076293945),
    @(10, 69.34999847412109, 82.44999694824219, 86.15000152587891, 67.15000152587891),
    @(12, 72.44999694824219, 69.80000305175781, 73.34999847412109, 67.05000305175781),
    @(13, 57.04999923706055, 61.04999923706055, 63.20000076293945, 56.95000076293945),
    @(14, 78.80000305175781, 82.25, 89.94999694824219, 75.65000152587891),
    @(15, 99.5, 95, 112.4000015258789, 94.75),
    @(16, 120.5299987792969, 97.34999847412109, 121.4499969482422, 87.93000030517578),
    @(17, 88, 109.3499984741211, 110.6500015258789, 84.83999633789062),
    @(18, 120.9499969482422, 134.1600036621094, 135.0700073242188, 115.3399963378906),
    @(19, 144.1199951171875, 148.5200042724609, 155.75, 139.1600036621094),
    @(20, 117.5800018310547, 122.0699996948242, 127.120002746582, 111.5999984741211),
    @(21, 123.2099990844727, 142.6900024414062, 145.3699951171875, 123.2099990844727),
    @(22, 97.12999725341795, 130.8099975585938, 137.3699951171875, 89.20999908447266),
    @(23, 258.4100036621094, 290.4800109863281, 299.8299865722656, 254.6799926757812),
    @(24, 259, 247.3200073242188, 298.739013671875, 244.8099975585937),
    @(25, 253.3200073242188, 247.0500030517578, 272.1499938964844, 236.8999938964844),
    @(26, 289.1499938964844, 317.8800048828125, 329, 283.8900146484375),
    @(27, 291.25, 298.6400146484375, 309, 272),
    @(28, 195.8999938964844, 185.9600067138672, 202.7050018310547, 171.3699951171875),
    @(29, 158.2100067138672, 131.3699951171875, 163.3899993896484, 113.4300003051758),
    @(30, 105.5, 75.45999908447266, 111.3499984741211, 73.82099914550781),
    @(31, 66.25, 59.33000183105469, 72.77999877929688, 56.16999816894531),
    @(32, 78.23999786376953, 84.09999847412109, 87.44000244140625, 67.37000274658203),
    @(33, 77.56999969482422, 86.98000335693359, 89.51000213623047, 70.27999877929688),
    @(34, 98.31999969482422, 87.23000335693359, 101.5500030517578, 84.12999725341797),
    @(35, 78.81999969482422, 94.31999969482422, 94.53600311279295, 75.08999633789062),
    @(36, 91.69000244140624, 79.90000152587891, 92.63300323486328, 77.51999664306641),
    @(37, 121.1999969482422, 126.879997253418, 133.0899963378906, 114.1800003051758),
    @(38, 137.6699981689453, 118.870002746582, 139, 117.7669982910156),
    @(39, 159.3099975585938, 155.9250030517578, 169.5800018310547, 150.6999969482422),
    @(40, 165.25, 167.1199951171875, 179, 155.7200012207031),
    @(41, 214.1799926757812, 238.8899993896484, 247.1100006103516, 208),
    @(42, 163.3800048828125, 169.5899963378906, 171.1100006103516, 142.3600006103516),
    @(43, 158.4600067138672, 136.0299987792969, 167.5200042724609, 135.3200073242188),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $ws.Cells.Item($r, 8).Value = 55686229
    $ws.Cells.Item($r, 9).Value = "WIX"
}
